# Updated remaining queries for C3DC
# Applies the SQL join-column renames (std.id/prt.id -> std.study_id/prt.participant_id,
# and the matching quoted "study.id"/"participant.id" -> "study.study_id"/"participant.participant_id")
# across every query cell (C2 and B2:B7) on Sheet1, plus the accompanying formatting
# clean-up (de-duplicated font list / column C width / B7 style index) that Excel produced
# when the file was re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix up the LEFT JOIN predicates in every SQL query cell (C2 and B2:B7).
# ---------------------------------------------------------------------------
$replacements = @(
    @('std.id = prt."study.id"',              'std.study_id = prt."study.study_id"'),
    @('prt.id = dgn."participant.id"',          'prt.participant_id = dgn."participant.participant_id"'),
    @('prt.id = trt."participant.id"',          'prt.participant_id = trt."participant.participant_id"'),
    @('prt.id = trr."participant.id"',          'prt.participant_id = trr."participant.participant_id"'),
    @('prt.id = srv."participant.id"',          'prt.participant_id = srv."participant.participant_id"'),
    @('std.id = rfs."study.id"',                'std.study_id = rfs."study.study_id"')
)

$targetCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $text = $cell.Value()
    foreach ($pair in $replacements) {
        $text = $text.Replace($pair[0], $pair[1])
    }
    $cell.Value = $text
}

# ---------------------------------------------------------------------------
# 2) Column C was widened (and is no longer an auto "best fit" width).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 69.6640625

# ---------------------------------------------------------------------------
# 3) B7 picked up a redundant duplicate font during earlier edits; re-applying
#    the same wrap-text style that B2:B6 already use collapses the duplicate
#    font entry Excel had kept around (fonts count 6 -> 5, cellXfs 5 -> 4).
# ---------------------------------------------------------------------------
$ws.Range("B7").Style = $ws.Range("B6").Style
